$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Ashmandeep Kaur"

# Test plan rows: Preconditions (E), Method Inputs (F), Expected Result (G)
$ws.Range("E7").Value  = "None"
$ws.Range("F7").Value  = "account_number=12345, client_number=67890, balance=1000.00"
$ws.Range("G7").Value  = "BankAccount instance is created successfully with correct values."

$ws.Range("E8").Value  = "None"
$ws.Range("F8").Value  = "balance='invalid_balance'"
$ws.Range("G8").Value  = "balance is set to 0."

$ws.Range("E9").Value  = "None"
$ws.Range("F9").Value  = "account_number='abc'"
$ws.Range("G9").Value  = "Raises ValueError with message about invalid account number"

$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "client_number='xyz'"
$ws.Range("G10").Value = "Raises ValueError with message about invalid client number."

$ws.Range("E11").Value = "BankAccount instance exists"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "Returns the account number of the instance."

$ws.Range("E12").Value = "BankAccount instance exists"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "Returns the client number of the instance."

$ws.Range("E13").Value = "BankAccount instance exists"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Returns the balance of the instance."

$ws.Range("E14").Value = "BankAccount instance exists"
$ws.Range("F14").Value = "amount=500.00"
$ws.Range("G14").Value = "Balance is updated by adding 500 to current balance."

$ws.Range("E15").Value = "BankAccount instance exists"
$ws.Range("F15").Value = "amount=200.00"
$ws.Range("G15").Value = "Balance is updated by deducting 200 from current balance."

$ws.Range("E16").Value = "BankAccount instance exists"
$ws.Range("F16").Value = "amount='invalid'"
$ws.Range("G16").Value = "Balance remains unchanged; non-numeric input does not update balance."

$ws.Range("E17").Value = "BankAccount instance exists"
$ws.Range("F17").Value = "amount=300.00"
$ws.Range("G17").Value = "Balance is updated by adding 300 to the current balance."

$ws.Range("E18").Value = "BankAccount instance exists"
$ws.Range("F18").Value = "amount=100.00"
$ws.Range("G18").Value = "Raises ValueError with a message stating that deposit amount must be positive."

$ws.Range("E19").Value = "BankAccount instance exists"
$ws.Range("F19").Value = "amount=200.00"
$ws.Range("G19").Value = "Balance is updated by deducting 200 from the current balance."

$ws.Range("E20").Value = "BankAccount instance exists"
$ws.Range("F20").Value = "amount=50.00"
$ws.Range("G20").Value = "Raises ValueError with a message stating withdrawal amount must be positive."

$ws.Range("E21").Value = "BankAccount instance exists"
$ws.Range("F21").Value = "account=5000.00(more than current balance)"
$ws.Range("G21").Value = "Raises ValueError with a message stating withdrawal exceeds the account balance."

$ws.Range("E22").Value = "BankAccount instance exists"
$ws.Range("F22").Value = "None"
$ws.Range("G22").Value = "Returns account details as formatted string with balance in currency format."

# Match the final selection recorded in the workbook
$ws.Range("G22").Select()
